$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "21.651.88"
$ws.Range("E2").Value = "  -1.93%  "

# Row 3
$ws.Range("D3").Value = "1.532.49"
$ws.Range("E3").Value = "  -1.58%  "

# Row 4
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.18%  "

# Row 5
$ws.Range("E5").Value = "  +0.23%  "

# Row 6
$ws.Range("D6").Value = "'288.44"
$ws.Range("E6").Value = "  +0.42%  "

# Row 7
$ws.Range("D7").Value = "'0.3944"
$ws.Range("E7").Value = "  +2.73%  "

# Row 8
$ws.Range("D8").Value = "'0.3147"
$ws.Range("E8").Value = "  -3.09%  "

# Row 9
$ws.Range("D9").Value = "'42.32"
$ws.Range("E9").Value = "  +2.12%  "

# Row 10
$ws.Range("D10").Value = "'0.07139"
$ws.Range("E10").Value = "  -2.69%  "

# Row 11
$ws.Range("D11").Value = "'1.044"
$ws.Range("E11").Value = "  -7.44%  "

# Row 12
$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "  +0.31%  "

# Row 13
$ws.Range("D13").Value = "'5.641"
$ws.Range("E13").Value = "  -1.48%  "

# Row 14
$ws.Range("D14").Value = "'18.51"
$ws.Range("E14").Value = "  -4.93%  "

# Row 15
$ws.Range("D15").Value = "'6.577"
$ws.Range("E15").Value = "  -3.54%  "

# Row 16
$ws.Range("D16").Value = "1.540.87"
$ws.Range("E16").Value = "  -1.14%  "

# Row 17
$ws.Range("D17").Value = "'0.00001086"
$ws.Range("E17").Value = "  -0.53%  "

# Row 18
$ws.Range("D18").Value = "'0.06594"
$ws.Range("E18").Value = "  -0.37%  "

# Row 19
$ws.Range("D19").Value = "'83.22"
$ws.Range("E19").Value = "  -2.20%  "

# Row 20
$ws.Range("E20").Value = "  +0.19%  "

# Row 21
$ws.Range("D21").Value = "'6.100"
$ws.Range("E21").Value = "  -4.89%  "

# Row 22
$ws.Range("D22").Value = "'15.38"
$ws.Range("E22").Value = "  -3.86%  "

# Row 23
$ws.Range("E23").Value = "  -6.14%  "

# Row 24
$ws.Range("D24").Value = "'2.361"
$ws.Range("E24").Value = "  +1.27%  "

# Row 25
$ws.Range("D25").Value = "21.635.72"
$ws.Range("E25").Value = "  -2.06%  "

# Row 26
$ws.Range("D26").Value = "'2.337"
$ws.Range("E26").Value = "  -8.25%  "

# Row 27
$ws.Range("D27").Value = "'147.87"
$ws.Range("E27").Value = "  -0.83%  "

# Row 28
$ws.Range("D28").Value = "'18.33"
$ws.Range("E28").Value = "  -2.89%  "

# Row 29
$ws.Range("D29").Value = "'4.840"
$ws.Range("E29").Value = "  -0.30%  "

# Row 30
$ws.Range("D30").Value = "1.707.98"
$ws.Range("E30").Value = "  -1.33%  "

# Row 31
$ws.Range("D31").Value = "'116.82"
$ws.Range("E31").Value = "  -3.33%  "

# Row 32
$ws.Range("D32").Value = "'5.850"
$ws.Range("E32").Value = "  -1.10%  "

# Row 33
$ws.Range("D33").Value = "'0.9414"
$ws.Range("E33").Value = "  -15.16%  "

# Row 34
$ws.Range("D34").Value = "'0.08131"
$ws.Range("E34").Value = "  -0.47%  "

# Row 35
$ws.Range("D35").Value = "'8.499"
$ws.Range("E35").Value = "  -8.48%  "

# Row 36
$ws.Range("D36").Value = "'5.110"
$ws.Range("E36").Value = "  -2.64%  "

# Row 37
$ws.Range("D37").Value = "'0.05991"
$ws.Range("E37").Value = "  -3.55%  "

# Row 38
$ws.Range("D38").Value = "'0.02195"
$ws.Range("E38").Value = "  -4.46%  "

# Row 39
$ws.Range("D39").Value = "'1.440"
$ws.Range("E39").Value = "  -14.01%  "

# Row 40
$ws.Range("D40").Value = "'0.2013"
$ws.Range("E40").Value = "  -4.76%  "

# Row 41
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "'10.96"
$ws.Range("E41").Value = "  +0.32%  "

# Row 42
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.167"
$ws.Range("E42").Value = "  -4.69%  "

# Row 43
$ws.Range("B43").Value = "Frax"
$ws.Range("C43").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D43").Value = "'1.000"
$ws.Range("E43").Value = "  +0.20%  "

# Row 44
$ws.Range("D44").Value = "'0.5730"
$ws.Range("E44").Value = "  -3.71%  "

# Row 45
$ws.Range("D45").Value = "'12.93"
$ws.Range("E45").Value = "  -4.40%  "

# Row 46
$ws.Range("D46").Value = "'3.714"
$ws.Range("E46").Value = "  -0.16%  "

# Row 47
$ws.Range("D47").Value = "'0.5474"
$ws.Range("E47").Value = "  -4.90%  "

# Row 48
$ws.Range("D48").Value = "'1.165"
$ws.Range("E48").Value = "  +0.51%  "

# Row 49
$ws.Range("D49").Value = "'116.04"
$ws.Range("E49").Value = "  -2.94%  "

# Row 50
$ws.Range("D50").Value = "'1.859"
$ws.Range("E50").Value = "  -4.11%  "

# Row 51
$ws.Range("D51").Value = "'0.06682"
$ws.Range("E51").Value = "  -3.11%  "

